{"js": "// Fix three Serbian-language typos in the document:\n//  1) \"Tok dogadjaja\"  -> \"Tok doga\u0111aja\"   (both the ToC entry and the heading)\n//  2) \"Takodje\"        -> \"Tako\u0111e\"         (inside the \"Dokumentom...\" paragraph)\n//  3) \"zatevanih\"      -> \"zahtevanih\"     (only inside the Heading 3 \"2.2.1.4.a ...\" title,\n//                                           NOT inside its Table-of-Contents entry)\n\n// 1) \"Tok dogadjaja\" -> \"Tok doga\u0111aja\" (appears twice: ToC line + heading, both change)\nconst dogadjaja = context.document.body.search(\"Tok dogadjaja\", { matchCase: true });\ndogadjaja.load(\"items\");\nawait context.sync();\nfor (const r of dogadjaja.items) {\n  r.insertText(\"Tok doga\u0111aja\", \"Replace\");\n}\n\n// 2) \"Takodje\" -> \"Tako\u0111e\"\nconst takodje = context.document.body.search(\"Takodje\", { matchCase: true });\ntakodje.load(\"items\");\nawait context.sync();\nfor (const r of takodje.items) {\n  r.insertText(\"Tako\u0111e\", \"Replace\");\n}\n\n// 3) \"zatevanih\" -> \"zahtevanih\", only for the occurrence that lives in a Heading 3\n//    paragraph (the Table-of-Contents occurrence must stay untouched).\nconst zatevanih = context.document.body.search(\"zatevanih\", { matchCase: true });\nzatevanih.load(\"items\");\nawait context.sync();\n\nconst zatevanihParas = zatevanih.items.map((r) => r.paragraphs.getFirst());\nzatevanihParas.forEach((p) => p.load(\"style\"));\nawait context.sync();\n\nfor (let i = 0; i < zatevanih.items.length; i++) {\n  if (zatevanihParas[i].style === \"Heading 3\") {\n    zatevanih.items[i].insertText(\"zahtevanih\", \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix three Serbian-language typos in the document:\n#  1) \"Tok dogadjaja\"  -> \"Tok doga\u0111aja\"   (both the ToC entry and the heading)\n#  2) \"Takodje\"        -> \"Tako\u0111e\"         (inside the \"Dokumentom...\" paragraph)\n#  3) \"zatevanih\"      -> \"zahtevanih\"     (only inside the Heading 3 \"2.2.1.4.a ...\" title,\n#                                           NOT inside its Table-of-Contents entry)\n\n$d = $word.ActiveDocument\n\n# 1) \"Tok dogadjaja\" -> \"Tok doga\u0111aja\" (replace every occurrence: ToC line + heading)\n$find = $d.Content.Find\n$find.Text = \"Tok dogadjaja\"\n$find.Replacement.Text = \"Tok doga\u0111aja\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) \"Takodje\" -> \"Tako\u0111e\"\n$find2 = $d.Content.Find\n$find2.Text = \"Takodje\"\n$find2.Replacement.Text = \"Tako\u0111e\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# 3) \"zatevanih\" -> \"zahtevanih\", only inside the Heading 3 paragraph\n#    (the Table-of-Contents occurrence of the same word must stay untouched).\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Heading 3\" -and $p.Range.Text -like \"*zatevanih*\") {\n        $find3 = $p.Range.Find\n        $find3.Text = \"zatevanih\"\n        $find3.Replacement.Text = \"zahtevanih\"\n        $find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2) | Out-Null\n    }\n}\n"}
